# Natmi following Dr Hou advice
#
# The ligand/receptor cluster-pair table is recomputed: a third "sending"
# / "target" cluster group (ECs) is added alongside the existing FAPs/sCs
# groups, turning the 2x3 cluster-pair grid (rows 2-7) into a 3x3 grid
# (rows 2-10), with refreshed statistics for every pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A..T =
#  Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
#  Ligand-expressing cells, Ligand detection rate,
#  Ligand average expression value, Ligand total expression value,
#  Ligand derived specificity of average expression value,
#  Ligand derived specificity of total expression value,
#  Receptor-expressing cells, Receptor detection rate,
#  Receptor average expression value, Receptor total expression value,
#  Receptor derived specificity of average expression value,
#  Receptor derived specificity of total expression value,
#  Edge average expression weight, Edge total expression weight,
#  Edge average expression derived specificity,
#  Edge total expression derived specificity

$rows = @()
$rows += ,@("ECs",  "Tgfb3", "Tgfbr2", "ECs",  2, 0.6666666666666666, 0.6423476666666667, 1.927043,           0.01173234890143342, 0.01173234890143342, 3, 1, 30.84053233333333, 92.521597,         0.2792762306509579, 0.2792762306509579, 19.81034398307455, 178.293095847671,   0.003276566177874232, 0.003276566177874233)
$rows += ,@("ECs",  "Tgfb3", "Tgfbr2", "FAPs", 2, 0.6666666666666666, 0.6423476666666667, 1.927043,           0.01173234890143342, 0.01173234890143342, 3, 1, 69.52746833333333, 208.582405,        0.6296055163046042, 0.6296055163046042, 44.66080705315722, 401.947263478415,   0.007386751587552743, 0.007386751587552745)
$rows += ,@("ECs",  "Tgfb3", "Tgfbr2", "sCs",  2, 0.6666666666666666, 0.6423476666666667, 1.927043,           0.01173234890143342, 0.01173234890143342, 3, 1, 10.06220766666667, 30.186623,         0.09111825304443795, 0.09111825304443795, 6.463435616198777, 58.170920545789, 0.001069031136006444, 0.001069031136006444)
$rows += ,@("FAPs", "Tgfb3", "Tgfbr2", "ECs",  3, 1,                  22.27635266666667,  66.829058,          0.4068730304461968,  0.4068730304461968,  3, 1, 30.84053233333333, 92.521597,         0.2792762306509579, 0.2792762306509579, 687.0145746850695, 6183.131172165627, 0.1136299662965463,   0.1136299662965463)
$rows += ,@("FAPs", "Tgfb3", "Tgfbr2", "FAPs", 3, 1,                  22.27635266666667,  66.829058,          0.4068730304461968,  0.4068730304461968,  3, 1, 69.52746833333333, 208.582405,        0.6296055163046042, 0.6296055163046042, 1548.818404613832, 13939.36564152449, 0.2561695044044967,   0.2561695044044967)
$rows += ,@("FAPs", "Tgfb3", "Tgfbr2", "sCs",  3, 1,                  22.27635266666667,  66.829058,          0.4068730304461968,  0.4068730304461968,  3, 1, 10.06220766666667, 30.186623,         0.09111825304443795, 0.09111825304443795, 224.1492865879038, 2017.343579291134, 0.03707355974515386,  0.03707355974515386)
$rows += ,@("sCs",  "Tgfb3", "Tgfbr2", "ECs",  3, 1,                  31.831433,          95.494299,          0.5813946206523698,  0.5813946206523697,  3, 1, 30.84053233333333, 92.521597,         0.2792762306509579, 0.2792762306509579, 981.6983386528336, 8835.285047875503, 0.1623696981765374,   0.1623696981765374)
$rows += ,@("sCs",  "Tgfb3", "Tgfbr2", "FAPs", 3, 1,                  31.831433,          95.494299,          0.5813946206523698,  0.5813946206523697,  3, 1, 69.52746833333333, 208.582405,        0.6296055163046042, 0.6296055163046042, 2213.158949912121, 19918.43054920909, 0.3660492603125547,   0.3660492603125547)
$rows += ,@("sCs",  "Tgfb3", "Tgfbr2", "sCs",  3, 1,                  31.831433,          95.494299,          0.5813946206523698,  0.5813946206523697,  3, 1, 10.06220766666667, 30.186623,         0.09111825304443795, 0.09111825304443795, 320.2944891735863, 2882.650402562277, 0.05297566216327764,  0.05297566216327763)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $rowNum = $i + 2
    $rowData = $rows[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $colNum = $j + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $rowData[$j]
    }
}
